$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "57.512.52"
$ws.Range("E2").Value = "  -2.03%  "
# Row 3
$ws.Range("D3").Value = "2.433.50"
$ws.Range("E3").Value = "  -2.09%  "
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.14%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "515.07"
$ws.Range("E5").Value = "  -2.15%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.28"
$ws.Range("E6").Value = "  -2.29%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.19%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.550"
$ws.Range("E8").Value = "  -1.97%  "
# Row 9
$ws.Range("D9").Value = "2.445.94"
$ws.Range("E9").Value = "  -1.50%  "
# Row 10
$ws.Range("E10").Value = "  -0.41%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0958"
$ws.Range("E11").Value = "  -4.40%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.20"
$ws.Range("E12").Value = "  -2.88%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.334"
$ws.Range("E13").Value = "  -2.03%  "
# Row 14
$ws.Range("D14").Value = "2.861.81"
$ws.Range("E14").Value = "  -2.21%  "
# Row 15
$ws.Range("D15").Value = "57.425.87"
$ws.Range("E15").Value = "  -2.05%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.85"
$ws.Range("E16").Value = "  -2.65%  "
# Row 17
$ws.Range("E17").Value = "  -3.21%  "
# Row 18
$ws.Range("D18").Value = "2.435.12"
$ws.Range("E18").Value = "  -2.44%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.48"
$ws.Range("E19").Value = "  -3.72%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "315.79"
$ws.Range("E20").Value = "  -1.49%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.11"
$ws.Range("E21").Value = "  -2.14%  "
# Row 22
$ws.Range("E22").Value = "  +0.11%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.69"
$ws.Range("E23").Value = "  -2.13%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.05"
$ws.Range("E24").Value = "  -0.05%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.406"
$ws.Range("E25").Value = "  -1.84%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.990"
$ws.Range("E26").Value = "  -1.14%  "
# Row 27
$ws.Range("E27").Value = "  -2.72%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.22"
$ws.Range("E28").Value = "  -2.97%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "168.39"
$ws.Range("E29").Value = "  +1.21%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.29"
$ws.Range("E30").Value = "  -1.68%  "
# Row 31
$ws.Range("B31").Value = "PEPE"
$ws.Range("C31").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D31").Value = "0.0₃0723"
$ws.Range("E31").Value = "  -3.58%  "
# Row 32
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.67"
$ws.Range("E32").Value = "  -2.72%  "
# Row 33
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.17"
$ws.Range("E33").Value = "  +3.53%  "
# Row 34
$ws.Range("E34").Value = "  -0.05%  "
# Row 35
$ws.Range("E35").Value = "  -0.21%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.75"
$ws.Range("E36").Value = "  -2.67%  "
# Row 37
$ws.Range("E37").Value = "  -4.81%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.92"
$ws.Range("E38").Value = "  -1.14%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.20"
$ws.Range("E39").Value = "  -1.43%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.46"
$ws.Range("E40").Value = "  -2.48%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.782"
$ws.Range("E41").Value = "  -0.75%  "
# Row 42
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.02"
$ws.Range("E42").Value = "  +1.86%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "269.91"
$ws.Range("E43").Value = "  -2.04%  "
# Row 44
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.38"
$ws.Range("E44").Value = "  -4.10%  "
# Row 45
$ws.Range("E45").Value = "  -1.50%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "121.27"
$ws.Range("E46").Value = "  -4.77%  "
# Row 47
$ws.Range("E47").Value = "  -0.76%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0485"
$ws.Range("E48").Value = "  -1.66%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "17.13"
$ws.Range("E49").Value = "  -3.09%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0210"
$ws.Range("E50").Value = "  -2.20%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.56"
$ws.Range("E51").Value = "  -2.91%  "
